$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full contact list (row 1 = header, rows 2-24 = contacts)
$data = @(
    @('+5511977645543', 'Vinicius Simulacro'),
    @('+5513991892211', 'aloha bar'),
    @('+5561995672332', 'Porks'),
    @('+5513974226875', 'Beco Bar'),
    @('+5513998027988', 'Maleh'),
    @('+5513996000859', 'Resenhas'),
    @('+5513996672625', 'Barzin da praia'),
    @('+5513981840830', 'Avenisdas Bar e espeto'),
    @('+5513991328222', 'Confraria do forte'),
    @('+5513996935710', 'Quiosque 7'),
    @('+5513991342205', 'Vixe Bar'),
    @('+5513974019871', 'Matilde Bar'),
    @('+5513991037423', 'Capitão Bar'),
    @('+5513981356526', 'neco''s bar'),
    @('+5513997005002', 'wall street hamb'),
    @('+5513996009483', 'Nabarca sushi'),
    @('+5513991302000', 'Elo Adega e Tabaca'),
    @('+5513991705555', 'Baroni Bar'),
    @('+5513974098864', 'Pink Bar'),
    @('+5513988659250', 'Espetos Japa e Cia'),
    @('+5513997529370', 'Bar se7e'),
    @('+551335916484', 'Gaucho Mallet Gril'),
    @('+551333295300', 'Santo Canto')
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("A23").Select() | Out-Null
